$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 - "Ann Window Volume" for Panel B (E-mini Futures) / Emini
$ws.Range("D26").Value = 2136.557752341311
$ws.Range("E26").Value = 883.4819743009047
$ws.Range("F26").Value = 1415.153225806452
$ws.Range("G26").Value = 2099.241935483871
$ws.Range("H26").Value = 2698.580645161291
$ws.Range("I26").Value = 62
$ws.Range("J26").Value = 2064.061872025384
$ws.Range("K26").Value = 817.6011738785187
$ws.Range("L26").Value = 1451.331967213115
$ws.Range("M26").Value = 2061.55737704918
$ws.Range("N26").Value = 2642.516393442623
$ws.Range("O26").Value = 62
$ws.Range("P26").Value = 2052.031458277792
$ws.Range("Q26").Value = 800.9504474134461
$ws.Range("R26").Value = 1498.423553719008
$ws.Range("S26").Value = 1975.152892561984
$ws.Range("T26").Value = 2549.993801652892
$ws.Range("U26").Value = 62
$ws.Range("V26").Value = 2070.525883256529
$ws.Range("W26").Value = 764.2420525424247
$ws.Range("X26").Value = 1465.779761904762
$ws.Range("Y26").Value = 1954.283333333333
$ws.Range("Z26").Value = 2610.178571428572
$ws.Range("AA26").Value = 62
$ws.Range("AB26").Value = 800.3816593352884
$ws.Range("AC26").Value = 227.8441033921607
$ws.Range("AD26").Value = 652.3098484848485
$ws.Range("AE26").Value = 750.2329545454545
$ws.Range("AF26").Value = 981.1994318181818
$ws.Range("AG26").Value = 62

# Row 27 - "Diff (Ann - Non)" for Panel B (E-mini Futures) / Emini
$ws.Range("D27").Value = 441.556581685744
$ws.Range("J27").Value = 358.3273400317292
$ws.Range("P27").Value = 315.4530791788857
$ws.Range("V27").Value = 256.4359447004608
$ws.Range("AB27").Value = 65.97499389051806

# Row 28 - "# Obs" for Panel B (E-mini Futures) / Emini
$ws.Range("D28").Value = 62
$ws.Range("J28").Value = 62
$ws.Range("P28").Value = 62
$ws.Range("V28").Value = 62
$ws.Range("AB28").Value = 62
